$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.949795603752136
$ws.Range("B1").Value = 2.139410972595215
$ws.Range("C1").Value = 2.519808292388916
$ws.Range("D1").Value = 3.155150651931763
$ws.Range("E1").Value = 2.27548623085022
